$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.271.23"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "2.331.59"
$ws.Range("E3").Value = "  -0.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "533.93"
$ws.Range("E5").Value = "  +2.92%  "

# Row 6
$ws.Range("D6").Value = "132.83"
$ws.Range("E6").Value = "  -2.26%  "

# Row 7
$ws.Range("E7").Value = "  -0.48%  "

# Row 8
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").Value = "2.359.75"
$ws.Range("E9").Value = "  +0.38%  "

# Row 10
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  -1.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("D12").Value = "5.32"
$ws.Range("E12").Value = "  -1.81%  "

# Row 13
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "23.63"
$ws.Range("E14").Value = "  -1.43%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.750.39"
$ws.Range("E15").Value = "  -0.29%  "

# Row 16
$ws.Range("D16").Value = "57.294.16"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17
$ws.Range("E17").Value = "  -1.24%  "

# Row 18
$ws.Range("D18").Value = "2.346.96"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("D19").Value = "339.41"
$ws.Range("E19").Value = "  +3.74%  "

# Row 20
$ws.Range("D20").Value = "10.46"
$ws.Range("E20").Value = "  -1.47%  "

# Row 21
$ws.Range("E21").Value = "  +2.30%  "

# Row 22
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  -1.49%  "

# Row 23
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").Value = "61.83"
$ws.Range("E24").Value = "  +1.08%  "

# Row 25
$ws.Range("D25").Value = "8.91"
$ws.Range("E25").Value = "  +10.82%  "

# Row 26
$ws.Range("E26").Value = "  -0.22%  "

# Row 27
$ws.Range("D27").Value = "0.992"
$ws.Range("E27").Value = "  -0.37%  "

# Row 28
$ws.Range("E28").Value = "  +3.39%  "

# Row 29
$ws.Range("D29").Value = "169.83"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").Value = "1.73"
$ws.Range("E30").Value = "  +1.60%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0729"
$ws.Range("E31").Value = "  -1.85%  "

# Row 32
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
$ws.Range("D33").Value = "18.53"
$ws.Range("E33").Value = "  -0.21%  "

# Row 34
$ws.Range("D34").Value = "0.999"

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -0.58%  "

# Row 37
$ws.Range("E37").Value = "  +0.22%  "

# Row 38
$ws.Range("D38").Value = "0.913"
$ws.Range("E38").Value = "  -0.24%  "

# Row 39
$ws.Range("E39").Value = "  +1.21%  "

# Row 40
$ws.Range("E40").Value = "  +1.32%  "

# Row 41
$ws.Range("D41").Value = "148.14"
$ws.Range("E41").Value = "  -1.10%  "

# Row 42
$ws.Range("D42").Value = "0.377"
$ws.Range("E42").Value = "  -1.70%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "286.84"
$ws.Range("E43").Value = "  +2.31%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.60"
$ws.Range("E44").Value = "  -1.40%  "

# Row 45
$ws.Range("D45").Value = "5.14"
$ws.Range("E45").Value = "  -1.26%  "

# Row 46
$ws.Range("E46").Value = "  -0.30%  "

# Row 47
$ws.Range("D47").Value = "0.0504"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "18.82"
$ws.Range("E48").Value = "  +4.35%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.562"
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").Value = "0.0218"
$ws.Range("E50").Value = "  -0.80%  "

# Row 51
$ws.Range("D51").Value = "17.38"
$ws.Range("E51").Value = "  -0.04%  "
